$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 44.51852
$ws.Range("I5").Value = 34.782608
$ws.Range("J5").Value = 100.5
$ws.Range("K5").Value = 34.782608
$ws.Range("L5").Value = 100.5
$ws.Range("M5").Value = 80.21739199999999
$ws.Range("N5").Value = -330.5
$ws.Range("H32").Value = 1528.2941
$ws.Range("I32").Value = 1620.0714
$ws.Range("K32").Value = 1620.0714
$ws.Range("M32").Value = -1294.0714
$ws.Range("H40").Value = 3042139
$ws.Range("I40").Value = 17142.715
$ws.Range("J40").Value = 8335882.5
$ws.Range("K40").Value = 17142.715
$ws.Range("L40").Value = 8335882.5
$ws.Range("M40").Value = -16967.715
$ws.Range("N40").Value = -8336232.5
$ws.Range("H64").Value = 6000
$ws.Range("J64").Value = 6000
$ws.Range("L64").Value = 6000
$ws.Range("N64").Value = -6496
$ws.Range("H67").Value = 6000
$ws.Range("J67").Value = 6000
$ws.Range("L67").Value = 6000
$ws.Range("N67").Value = -7716
$ws.Range("H74").Value = 55565776
$ws.Range("I74").Value = 71437144
$ws.Range("J74").Value = 16000
$ws.Range("K74").Value = 71437144
$ws.Range("L74").Value = 16000
$ws.Range("M74").Value = -71436208
$ws.Range("N74").Value = -17872
$ws.Range("H77").Value = 55565776
$ws.Range("I77").Value = 71437144
$ws.Range("J77").Value = 16000
$ws.Range("K77").Value = 357185720
$ws.Range("L77").Value = 80000
$ws.Range("M77").Value = -357181040
$ws.Range("N77").Value = -89360
$ws.Range("H80").Value = 71714.14
$ws.Range("J80").Value = 80299.8
$ws.Range("L80").Value = 240899.4
$ws.Range("N80").Value = -242895.4
$ws.Range("H83").Value = 71714.14
$ws.Range("J83").Value = 80299.8
$ws.Range("L83").Value = 722698.2000000001
$ws.Range("N83").Value = -732682.2000000001
$ws.Range("H98").Value = 6098.3887
$ws.Range("I98").Value = 5426.5713
$ws.Range("K98").Value = 5426.5713
$ws.Range("M98").Value = -3928.5713
$ws.Range("H106").Value = 3528.8333
$ws.Range("I106").Value = 3445.6365
$ws.Range("K106").Value = 3445.6365
$ws.Range("M106").Value = -2814.6365
$ws.Range("H113").Value = 146466640
$ws.Range("J113").Value = 166668160
$ws.Range("L113").Value = 166668160
$ws.Range("N113").Value = -166674668
$ws.Range("H121").Value = 3596.3572
$ws.Range("J121").Value = 3596.3572
$ws.Range("L121").Value = 10789.0716
$ws.Range("N121").Value = -14283.0716
$ws.Range("H122").Value = 6098.3887
$ws.Range("I122").Value = 5426.5713
$ws.Range("K122").Value = 16279.7139
$ws.Range("M122").Value = -13829.7139
$ws.Range("H138").Value = 2505414.8
$ws.Range("I138").Value = 2634.0667
$ws.Range("J138").Value = 4007083
$ws.Range("K138").Value = 7902.2001
$ws.Range("L138").Value = 12021249
$ws.Range("M138").Value = -2762.2001
$ws.Range("N138").Value = -12031529

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 7105
$ws.Range("I10").Value = 7105
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 7105
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = ""   # was -3945, cell removed
$ws.Range("M10").Value = -6935
$ws.Range("H32").Value = 3151.4814
$ws.Range("I32").Value = 3156.625
$ws.Range("K32").Value = 3156.625
$ws.Range("M32").Value = -2869.625
$ws.Range("H45").Value = 5999
$ws.Range("I45").Value = 3001
$ws.Range("J45").Value = 6598.6
$ws.Range("K45").Value = 3001
$ws.Range("L45").Value = 6598.6
$ws.Range("M45").Value = -2624
$ws.Range("N45").Value = -7352.6
$ws.Range("H61").Value = 4743.956
$ws.Range("I61").Value = 2604.2808
$ws.Range("K61").Value = 2604.2808
$ws.Range("M61").Value = -2392.2808
$ws.Range("H97").Value = 3796000.5
$ws.Range("I97").Value = 1403.0714
$ws.Range("J97").Value = 10436546
$ws.Range("K97").Value = 1403.0714
$ws.Range("L97").Value = 10436546
$ws.Range("M97").Value = -907.0714
$ws.Range("N97").Value = -10437538
$ws.Range("H122").Value = 13236.909
$ws.Range("I122").Value = 13695.368
$ws.Range("J122").Value = 10333.333
$ws.Range("K122").Value = 41086.104
$ws.Range("L122").Value = 30999.999
$ws.Range("M122").Value = -38636.104
$ws.Range("N122").Value = -35899.999
$ws.Range("H132").Value = 825732.9
$ws.Range("I132").Value = 1285005.2
$ws.Range("K132").Value = 3855015.6
$ws.Range("M132").Value = -3852485.6
$ws.Range("H136").Value = 4743.956
$ws.Range("I136").Value = 2604.2808
$ws.Range("K136").Value = 7812.8424
$ws.Range("M136").Value = -5262.8424

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 79900
$ws.Range("J74").Value = 79900
$ws.Range("L74").Value = 79900
$ws.Range("N74").Value = -81772
$ws.Range("H77").Value = 79900
$ws.Range("J77").Value = 79900
$ws.Range("L77").Value = 239700
$ws.Range("N77").Value = -249060
$ws.Range("H86").Value = 38502428
$ws.Range("J86").Value = 76926200
$ws.Range("L86").Value = 76926200
$ws.Range("N86").Value = -76928446
$ws.Range("H89").Value = 38502428
$ws.Range("J89").Value = 76926200
$ws.Range("L89").Value = 384631000
$ws.Range("N89").Value = -384642232
$ws.Range("H94").Value = 83337430
$ws.Range("J94").Value = 4837.875
$ws.Range("L94").Value = 4837.875
$ws.Range("N94").Value = -5739.875

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7240.109
$ws.Range("I31").Value = 3371.6956
$ws.Range("J31").Value = 11108.521
$ws.Range("K31").Value = 3371.6956
$ws.Range("L31").Value = 11108.521
$ws.Range("M31").Value = -3076.6956
$ws.Range("N31").Value = -11698.521
$ws.Range("H34").Value = 7240.109
$ws.Range("I34").Value = 3371.6956
$ws.Range("J34").Value = 11108.521
$ws.Range("K34").Value = 3371.6956
$ws.Range("L34").Value = 11108.521
$ws.Range("M34").Value = -3169.6956
$ws.Range("N34").Value = -11512.521
$ws.Range("H132").Value = 4927.4473
$ws.Range("I132").Value = 2532.6
$ws.Range("K132").Value = 7597.799999999999
$ws.Range("M132").Value = -5067.799999999999

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3317.5483
$ws.Range("J68").Value = 3543.8
$ws.Range("L68").Value = 10631.4
$ws.Range("N68").Value = -12253.4
$ws.Range("H71").Value = 3317.5483
$ws.Range("J71").Value = 3543.8
$ws.Range("L71").Value = 31894.2
$ws.Range("N71").Value = -40006.2
$ws.Range("H92").Value = 5918295
$ws.Range("J92").Value = 9616521
$ws.Range("L92").Value = 28849563
$ws.Range("N92").Value = -28852059
$ws.Range("H132").Value = 8092.2856
$ws.Range("J132").Value = 9138.684999999999
$ws.Range("L132").Value = 82248.16499999999
$ws.Range("N132").Value = -87308.16499999999

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = ""   # was -834, cell removed
$ws.Range("H11").Value = 552625
$ws.Range("I11").Value = 502250
$ws.Range("J11").Value = 603000
$ws.Range("K11").Value = 502250
$ws.Range("L11").Value = 603000
$ws.Range("M11").Value = -502111
$ws.Range("N11").Value = -603278
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -54900

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 10000
$ws.Range("I3").Value = 10000
$ws.Range("K3").Value = 10000
$ws.Range("M3").Value = -9888
$ws.Range("H15").Value = 10000
$ws.Range("I15").Value = 10000
$ws.Range("K15").Value = 10000
$ws.Range("M15").Value = -9830
$ws.Range("H40").Value = 4959.346
$ws.Range("I40").Value = 3643.4707
$ws.Range("K40").Value = 3643.4707
$ws.Range("M40").Value = -3507.4707
$ws.Range("H62").Value = 43666.668
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 43666.668
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 43666.668
$ws.Range("M62").Value = ""   # was -29376, cell removed
$ws.Range("N62").Value = -44914.668
$ws.Range("H65").Value = 43666.668
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 43666.668
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 131000.004
$ws.Range("M65").Value = ""   # was -86880, cell removed
$ws.Range("N65").Value = -137240.004
$ws.Range("H93").Value = 5324.75
$ws.Range("I93").Value = 5052.8667
$ws.Range("J93").Value = 5777.8887
$ws.Range("K93").Value = 5052.8667
$ws.Range("L93").Value = 5777.8887
$ws.Range("M93").Value = -3804.8667
$ws.Range("N93").Value = -8273.8887
$ws.Range("H136").Value = 10571.714
$ws.Range("I136").Value = 5084
$ws.Range("J136").Value = 17888.666
$ws.Range("K136").Value = 15252
$ws.Range("L136").Value = 53665.99800000001
$ws.Range("M136").Value = -12702
$ws.Range("N136").Value = -58765.99800000001
$ws.Range("H141").Value = 92911.8
$ws.Range("J141").Value = 92911.8
$ws.Range("L141").Value = 92911.8
$ws.Range("N141").Value = -103271.8

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 30000444
$ws.Range("I9").Value = 888
$ws.Range("K9").Value = 888
$ws.Range("M9").Value = -748
$ws.Range("H14").Value = 58826028
$ws.Range("I14").Value = 142859650
$ws.Range("K14").Value = 142859650
$ws.Range("M14").Value = -142859482
$ws.Range("H107").Value = 19609186
$ws.Range("I107").Value = 1451.4286
$ws.Range("J107").Value = 111111944
$ws.Range("K107").Value = 4354.2858
$ws.Range("L107").Value = 333335832
$ws.Range("M107").Value = -2434.2858
$ws.Range("N107").Value = -333339672

